$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$colLetters = @('B','D','E','F','G','H','J','L','M','N','O')

$data = @(
    @(2, 1.095202203248022, 0.3417225973513141, 0.3266785670077468, 1.063701366722306, 0.4948235029918493, 0.6208848012654329, 0.4015298274453798, 0.2246413329668542, 0.2549501717827383, 1.769497395984629, 2.193774985900433),
    @(3, 1.050540982382017, 0.3440791985611105, 0.3274015146307256, 1.065355026220061, 0.4886478917314037, 0.6217177850670765, 0.3997543040480593, 0.2016301905595128, 0.2387934273924301, 1.768414486219925, 2.182139606050981),
    @(4, 1.023408645877481, 0.3456134531433843, 0.3279113982133932, 1.066969736196498, 0.4852034929766944, 0.6225247511646046, 0.3987266893061872, 0.1874778665555255, 0.2289085101798136, 1.768376323061432, 2.176340052047209),
    @(5, 1.012425951366907, 0.3462606697953605, 0.32813585204544, 1.067778524063066, 0.4838872161587346, 0.6229279833317776, 0.3983238137255398, 0.1817052837766653, 0.2248895779670406, 1.768518965458583, 2.17431490309653),
    @(6, 1.010606781444949, 0.3463694694273487, 0.3281741316891544, 1.067921931421829, 0.4836739239781451, 0.6229994348207128, 0.398257880523083, 0.1807464392080931, 0.2242228065130902, 1.768552223216304, 2.173999059094541),
    @(7, 1.023260228564993, 0.3456220926142368, 0.3279143576671828, 1.066980033199371, 0.4851853875914287, 0.6225298880009547, 0.3987211914600479, 0.1874000366886861, 0.2288542715263517, 1.768377605472708, 2.176311370622926),
    @(8, 1.079743486471784, 0.3425170633165497, 0.3269141896459935, 1.064147190172605, 0.4926220198724991, 0.6211107265659734, 0.4009047329633333, 0.2167122296682464, 0.2493721876606472, 1.768994277366403, 2.189484083185619),
    @(9, 1.192765660987618, 0.3371186628569807, 0.3254730966148642, 1.06334586293498, 0.5099648123708818, 0.6206696552079762, 0.4056766715089637, 0.2739891466851248, 0.2898750828027801, 1.775152314202799, 2.225984782991389),
    @(10, 1.277133798135992, 0.3335705301245921, 0.3247270274536866, 1.065654644021024, 0.5243948421870641, 0.6217703461867785, 0.4094732642748511, 0.3159249837692357, 0.3197804095436325, 1.782663669819698, 2.259312954436155),
    @(11, 1.315794765838575, 0.3320465535065185, 0.3244545885188117, 1.067334128420541, 0.5313274285387735, 0.6225798869067916, 0.4112618681234324, 0.3349668721132275, 0.3334141999206395, 1.786723341002926, 2.275890514696385),
    @(12, 1.330474174919857, 0.331482371066588, 0.3243609727818004, 1.068060549430484, 0.5340056418504986, 0.6229307856327608, 0.4119478624374224, 0.3421720734468181, 0.3385809137753952, 1.78835251903854, 2.282371693129534),
    @(13, 1.327310969961275, 0.3316033041108732, 0.3243807110397086, 1.067900080165018, 0.5334264832933826, 0.62285324219404, 0.4117997373408571, 0.3406205595418044, 0.3374680023650498, 1.787997569256518, 2.280966800534713),
    @(14, 1.317001667737202, 0.3319998792808483, 0.3244466957174126, 1.067392079276743, 0.5315467045113849, 0.622607867232162, 0.4113181320475761, 0.3355597625671862, 0.3338391927166313, 1.786855536421484, 2.276419644744351),
    @(15, 1.310692011934464, 0.3322444736827022, 0.3244883547894517, 1.067092690305643, 0.530402187976776, 0.6224633408092473, 0.4110242617983886, 0.3324591401490693, 0.3316169378313631, 1.786167956036437, 2.273660895456004),
    @(16, 1.274612779838492, 0.3336719349378785, 0.3247461726622127, 1.065557543428049, 0.5239491952475674, 0.621723649602302, 0.4093575992153902, 0.314679803710959, 0.3188899746898457, 1.782411246429888, 2.258258068468962),
    @(17, 1.252550623665059, 0.3345706810749083, 0.3249214288977207, 1.06477690593956, 0.5200848461827263, 0.6213489249025628, 0.4083507992850741, 0.303763435259242, 0.3110897342059218, 1.780270846361972, 2.249171666606458),
    @(18, 1.239887616898898, 0.3350960978501201, 0.3250285380151361, 1.064387130599002, 0.517896838458924, 0.6211624657606336, 0.4077775139956046, 0.2974813719853557, 0.3066060609557368, 1.779100273152366, 2.24407873154388),
    @(19, 1.235604737097702, 0.3352754529143698, 0.325065888785609, 1.064265332887949, 0.5171619691142837, 0.6211043291627618, 0.4075844102904824, 0.2953538299472882, 0.3050884619057541, 1.778714349942305, 2.24237725139804),
    @(20, 1.254896437029913, 0.3344741304740606, 0.3249021204431912, 1.064853876696702, 0.5204926249334108, 0.6213858065608804, 0.408457375778319, 0.3049258421966101, 0.3119197941336012, 1.78049243518096, 2.250125131600953),
    @(21, 1.320028701083004, 0.3318830453636195, 0.3244270558161908, 1.067538837371899, 0.532097402627528, 0.6226787368733966, 0.4114593566308002, 0.3370463960780796, 0.3349049594234117, 1.787188490236687, 2.277749729053312),
    @(22, 1.362825176398189, 0.330264878273097, 0.3241722087107277, 1.069820728042728, 0.5399906953990268, 0.623782176930888, 0.4134718975655431, 0.3580064509978627, 0.3499496118144592, 1.792099978061501, 2.2969907994501),
    @(23, 1.339963334163428, 0.3311216516346103, 0.3243031600099755, 1.068554623554107, 0.5357496216675628, 0.6231696253023813, 0.4123931908505938, 0.3468228298254985, 0.3419180621701301, 1.789429832233651, 2.286612907350843),
    @(24, 1.253835829864556, 0.3345177538512187, 0.3249108300044163, 1.064818894360442, 0.5203081633209479, 0.6213690421017333, 0.4084091752856338, 0.3044003369741404, 0.3115445215537562, 1.780392068028476, 2.249693662145404),
    @(25, 1.161953146054714, 0.3385054461836128, 0.3258077167666507, 1.063053763786868, 0.5049770728293765, 0.6205385637039598, 0.4043341130683373, 0.2585182622153042, 0.2788909233220025, 1.772959716444518, 2.214967447626577)
)

foreach ($row in $data) {
    $r = [int]$row[0]
    for ($i = 0; $i -lt $colLetters.Count; $i++) {
        $addr = "$($colLetters[$i])$r"
        $ws.Range($addr).Value = [double]$row[$i + 1]
    }
}
